# Fruta / hortaliza, semanal
#
# Inserts 4 new weekly price-report rows for "Cebollín" (row 25..28),
# pushing the existing rows 25..48 down to 29..52, and extends the
# sheet's used range from A1:R48 to A1:R52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 4 new rows by inserting above the current row 25.
$ws.Range("A25:A28").EntireRow.Insert()

# --- New row 25: Primera, 300 kg, $600-$700 ($650 promedio) ---
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 44895
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = 100112037
$ws.Range("G25").Value = "Cebollín"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 600
$ws.Range("L25").Value = 700
$ws.Range("M25").Value = 650
$ws.Range("N25").Value = "$/paquete 6 unidades"
$ws.Range("O25").Value = "Provincia de Diguillín"
$ws.Range("P25").Value = 108
$ws.Range("Q25").Value = 6
$ws.Range("R25").Value = "Hortaliza"

# --- New row 26: Segunda, 300 kg, $500-$500 ($500 promedio) ---
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 44895
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 100112037
$ws.Range("G26").Value = "Cebollín"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 500
$ws.Range("L26").Value = 500
$ws.Range("M26").Value = 500
$ws.Range("N26").Value = "$/paquete 6 unidades"
$ws.Range("O26").Value = "Provincia de Diguillín"
$ws.Range("P26").Value = 83
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = "Hortaliza"

# --- New row 27: Primera, 400 kg, $600-$700 ($650 promedio) ---
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 44895
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 100112037
$ws.Range("G27").Value = "Cebollín"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 700
$ws.Range("M27").Value = 650
$ws.Range("N27").Value = "$/paquete 6 unidades"
$ws.Range("O27").Value = "Provincia de Diguillín"
$ws.Range("P27").Value = 108
$ws.Range("Q27").Value = 6
$ws.Range("R27").Value = "Hortaliza"

# --- New row 28: Segunda, 300 kg, $500-$500 ($500 promedio) ---
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44895
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112037
$ws.Range("G28").Value = "Cebollín"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Segunda"
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 500
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 500
$ws.Range("N28").Value = "$/paquete 6 unidades"
$ws.Range("O28").Value = "Provincia de Diguillín"
$ws.Range("P28").Value = 83
$ws.Range("Q28").Value = 6
$ws.Range("R28").Value = "Hortaliza"

Write-Output "Inserted 4 rows; dimension now $($ws.UsedRange.Rows.Count) rows"
